# Atualizado por script em 27-11-2023 08:45
#
# This script applies two kinds of changes to the Liga 1 Indonesia 2023-2024
# results sheet:
#   1. A handful of match rows had gotten interleaved out of chronological
#      order; for each affected pair of adjacent rows, the match details
#      (columns F..V: home/away teams, scores, odds, odds timestamps, url)
#      are swapped while the row index/metadata columns (A..E) stay put.
#   2. A new match result (Arema FC vs Persik Kediri) is appended as row 176.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the match-detail columns (F..V = column 6..22) between two rows,
# cell by cell (single-cell Value2 reads/writes are the reliable primitive
# here; multi-cell Range array assignment is not).
function Swap-MatchRows {
    param($RowA, $RowB)

    for ($col = 6; $col -le 22; $col++) {
        $a = $ws.Cells.Item($RowA, $col).Value2
        $b = $ws.Cells.Item($RowB, $col).Value2
        $ws.Cells.Item($RowA, $col).Value2 = $b
        $ws.Cells.Item($RowB, $col).Value2 = $a
    }
}

# Row pairs whose F:V contents (home/away teams, goals, odds, timestamps, url)
# need to be swapped with each other.
$pairs = @(
    @(17, 18),
    @(19, 20),
    @(42, 43),
    @(44, 45),
    @(82, 83),
    @(84, 85),
    @(93, 94),
    @(103, 104),
    @(107, 108),
    @(134, 135),
    @(148, 149),
    @(162, 163)
)

foreach ($pair in $pairs) {
    Swap-MatchRows $pair[0] $pair[1]
}

# Append the new match row (176) at the end of the sheet.
$newRow = 176
$prevRow = $newRow - 1

$ws.Cells.Item($newRow, 1).Value2 = 175
$ws.Cells.Item($newRow, 2).Value2 = "indonesia"
$ws.Cells.Item($newRow, 3).Value2 = "liga-1"
$ws.Cells.Item($newRow, 4).Value2 = "2023-2024"
$ws.Cells.Item($newRow, 5).Value2 = 45257.375

# Column A (Indice) and column E (data_partida) carry their own cell styles
# (bold/centered/bordered index column, and a date-time number format,
# respectively) throughout the sheet — copy those two cell formats down from
# the row above so the appended row matches the rest of the table.
$ws.Cells.Item($prevRow, 1).Copy() | Out-Null
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122) | Out-Null # xlPasteFormats
$ws.Cells.Item($prevRow, 5).Copy() | Out-Null
$ws.Cells.Item($newRow, 5).PasteSpecial(-4122) | Out-Null # xlPasteFormats

$ws.Cells.Item($newRow, 6).Value2 = "Arema FC"
$ws.Cells.Item($newRow, 7).Value2 = 0
$ws.Cells.Item($newRow, 8).Value2 = "Persik Kediri"
$ws.Cells.Item($newRow, 9).Value2 = 1
$ws.Cells.Item($newRow, 10).Value2 = 2.55
$ws.Cells.Item($newRow, 11).Value2 = "25/11/2023 21:12"
$ws.Cells.Item($newRow, 12).Value2 = 2.3
$ws.Cells.Item($newRow, 13).Value2 = "27/11/2023 07:56"
$ws.Cells.Item($newRow, 14).Value2 = 3.19
$ws.Cells.Item($newRow, 15).Value2 = "25/11/2023 21:12"
$ws.Cells.Item($newRow, 16).Value2 = 3.69
$ws.Cells.Item($newRow, 17).Value2 = "27/11/2023 07:22"
$ws.Cells.Item($newRow, 18).Value2 = 2.53
$ws.Cells.Item($newRow, 19).Value2 = "25/11/2023 21:12"
$ws.Cells.Item($newRow, 20).Value2 = 2.84
$ws.Cells.Item($newRow, 21).Value2 = "27/11/2023 07:56"
$ws.Cells.Item($newRow, 22).Value2 = "https://www.betexplorer.com/football/indonesia/liga-1/arema-fc-persik-kediri/2sYprumi/"
